# Scheduled runner update: refresh currentAveragePrice / LevePrice / LeveProfit
# figures across the Leve profitability sheets (market data resync).
#
# For each touched row we set the refreshed H:N values (currentAveragePrice,
# currentAveragePriceNQ, currentAveragePriceHQ, LevePriceNQ, LevePriceHQ,
# LeveProfitNQ, LeveProfitHQ). Where a column's cell no longer carries any
# value after the refresh, ClearContents() removes it instead of writing 0,
# matching rows where the source feed now omits that figure entirely.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H57").Value = 54195
$ws.Range("I57").Value = 27000
$ws.Range("J57").Value = 63260
$ws.Range("K57").Value = 81000
$ws.Range("L57").Value = 189780
$ws.Range("M57").Value = -80501
$ws.Range("N57").Value = -190778
$ws.Range("H132").Value = 1605.4
$ws.Range("I132").Value = 1561.5555
$ws.Range("K132").Value = 4684.666499999999
$ws.Range("M132").Value = -2154.666499999999
$ws.Range("H135").Value = 51725092
$ws.Range("I135").Value = 26316460
$ws.Range("J135").Value = 100001496
$ws.Range("K135").Value = 236848140
$ws.Range("L135").Value = 900013464
$ws.Range("M135").Value = -236845605
$ws.Range("N135").Value = -900018534
$ws.Range("H136").Value = 60577.418
$ws.Range("J136").Value = 60577.418
$ws.Range("L136").Value = 60577.418
$ws.Range("N136").Value = -70777.41800000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H45").Value = 1644.1163
$ws.Range("I45").Value = 1596.2258
$ws.Range("K45").Value = 1596.2258
$ws.Range("M45").Value = -1219.2258
$ws.Range("H74").Value = 6197.3335
$ws.Range("I74").Value = 2716.2778
$ws.Range("J74").Value = 11418.917
$ws.Range("K74").Value = 2716.2778
$ws.Range("L74").Value = 11418.917
$ws.Range("M74").Value = -1842.2778
$ws.Range("N74").Value = -13166.917
$ws.Range("H77").Value = 6197.3335
$ws.Range("I77").Value = 2716.2778
$ws.Range("J77").Value = 11418.917
$ws.Range("K77").Value = 13581.389
$ws.Range("L77").Value = 57094.585
$ws.Range("M77").Value = -9213.388999999999
$ws.Range("N77").Value = -65830.58499999999
$ws.Range("H134").Value = 47880.832
$ws.Range("J134").Value = 47880.832
$ws.Range("L134").Value = 47880.832
$ws.Range("N134").Value = -58020.832

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 4205.9287
$ws.Range("I134").Value = 4022.1428
$ws.Range("J134").Value = 4389.7144
$ws.Range("K134").Value = 12066.4284
$ws.Range("L134").Value = 13169.1432
$ws.Range("M134").Value = -9531.428400000001
$ws.Range("N134").Value = -18239.1432

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2680.4546
$ws.Range("I68").Value = 905.93335
$ws.Range("J68").Value = 4159.222
$ws.Range("K68").Value = 2717.80005
$ws.Range("L68").Value = 12477.666
$ws.Range("M68").Value = -1906.80005
$ws.Range("N68").Value = -14099.666
$ws.Range("H71").Value = 2680.4546
$ws.Range("I71").Value = 905.93335
$ws.Range("J71").Value = 4159.222
$ws.Range("K71").Value = 8153.40015
$ws.Range("L71").Value = 37432.998
$ws.Range("M71").Value = -4097.40015
$ws.Range("N71").Value = -45544.998
$ws.Range("H113").Value = 338.23
$ws.Range("I113").Value = 378.51514
$ws.Range("J113").Value = 318.38806
$ws.Range("K113").Value = 1135.54542
$ws.Range("L113").Value = 955.16418
$ws.Range("M113").Value = 1034.45458
$ws.Range("N113").Value = -5295.16418
$ws.Range("H120").Value = 12088.333
$ws.Range("I120").Value = 12088.333
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 36264.999
$ws.Range("L120").Value = 0
$ws.Range("M120").Value = -31426.999
$ws.Range("N120").ClearContents()
$ws.Range("H122").Value = 720.4545000000001
$ws.Range("I122").Value = 217.85715
$ws.Range("J122").Value = 1600
$ws.Range("K122").Value = 1960.71435
$ws.Range("L122").Value = 14400
$ws.Range("M122").Value = 489.28565
$ws.Range("N122").Value = -19300
$ws.Range("H123").Value = 1215
$ws.Range("I123").Value = 1215
$ws.Range("K123").Value = 3645
$ws.Range("M123").Value = -1195
$ws.Range("H124").Value = 1330.75
$ws.Range("J124").Value = 3033
$ws.Range("L124").Value = 9099
$ws.Range("N124").Value = -18919
$ws.Range("H125").Value = 3943.3333
$ws.Range("I125").Value = 2298
$ws.Range("J125").Value = 6000
$ws.Range("K125").Value = 6894
$ws.Range("L125").Value = 18000
$ws.Range("M125").Value = -1974
$ws.Range("N125").Value = -27840
$ws.Range("H133").Value = 4828.5
$ws.Range("I133").Value = 3087.7778
$ws.Range("J133").Value = 5750.0586
$ws.Range("K133").Value = 9263.3334
$ws.Range("L133").Value = 17250.1758
$ws.Range("M133").Value = -4203.3334
$ws.Range("N133").Value = -27370.1758
$ws.Range("H134").Value = 5070.36
$ws.Range("I134").Value = 7176.6665
$ws.Range("J134").Value = 3885.5625
$ws.Range("K134").Value = 21529.9995
$ws.Range("L134").Value = 11656.6875
$ws.Range("M134").Value = -16459.9995
$ws.Range("N134").Value = -21796.6875
$ws.Range("H136").Value = 3493.8235
$ws.Range("I136").Value = 936
$ws.Range("J136").Value = 4889
$ws.Range("K136").Value = 2808
$ws.Range("L136").Value = 14667
$ws.Range("M136").Value = 2292
$ws.Range("N136").Value = -24867
$ws.Range("H137").Value = 168732.67
$ws.Range("I137").Value = 3099
$ws.Range("J137").Value = 500000
$ws.Range("K137").Value = 9297
$ws.Range("L137").Value = 1500000
$ws.Range("M137").Value = -4197
$ws.Range("N137").Value = -1510200
$ws.Range("H139").Value = 5033992
$ws.Range("I139").Value = 5420762
$ws.Range("K139").Value = 16262286
$ws.Range("M139").Value = -16257146
$ws.Range("H141").Value = 3875.0435
$ws.Range("I141").Value = 2752.4
$ws.Range("J141").Value = 5980
$ws.Range("K141").Value = 8257.200000000001
$ws.Range("L141").Value = 17940
$ws.Range("M141").Value = -3077.200000000001
$ws.Range("N141").Value = -28300

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 9780
$ws.Range("I132").Value = 2491
$ws.Range("J132").Value = 53514
$ws.Range("K132").Value = 7473
$ws.Range("L132").Value = 160542
$ws.Range("M132").Value = -4943
$ws.Range("N132").Value = -165602
$ws.Range("H134").Value = 38334.57
$ws.Range("J134").Value = 38334.57
$ws.Range("L134").Value = 115003.71
$ws.Range("N134").Value = -120073.71
$ws.Range("H135").Value = 61010
$ws.Range("J135").Value = 61010
$ws.Range("L135").Value = 61010
$ws.Range("N135").Value = -71150

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 26000
$ws.Range("J101").Value = 26000
$ws.Range("L101").Value = 26000
$ws.Range("N101").Value = -32490
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H139").Value = 57452.5
$ws.Range("J139").Value = 57452.5
$ws.Range("L139").Value = 57452.5
$ws.Range("N139").Value = -67732.5

